$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 532.6667
$ws.Range("I103").Value = 374
$ws.Range("J103").Value = 850
$ws.Range("K103").Value = 1122
$ws.Range("L103").Value = 2550
$ws.Range("M103").Value = -536
$ws.Range("N103").Value = -3722
$ws.Range("H107").Value = 249.88889
$ws.Range("I107").Value = 241.26666
$ws.Range("J107").Value = 293
$ws.Range("K107").Value = 241.26666
$ws.Range("L107").Value = 293
$ws.Range("M107").Value = 1678.73334
$ws.Range("N107").Value = -4133
$ws.Range("H116").Value = 2444.4443
$ws.Range("I116").Value = 2142.8572
$ws.Range("K116").Value = 2142.8572
$ws.Range("M116").Value = 1299.1428
$ws.Range("H125").Value = 41667464
$ws.Range("I125").Value = 76923620
$ws.Range("J125").Value = 1097.2727
$ws.Range("K125").Value = 692312580
$ws.Range("L125").Value = 9875.454299999999
$ws.Range("M125").Value = -692310120
$ws.Range("N125").Value = -14795.4543
$ws.Range("H141").Value = 3059.0908
$ws.Range("I141").Value = 1838.0769
$ws.Range("J141").Value = 4822.778
$ws.Range("K141").Value = 5514.2307
$ws.Range("L141").Value = 14468.334
$ws.Range("M141").Value = -334.2307000000001
$ws.Range("N141").Value = -24828.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1250.6
$ws.Range("I2").Value = 945
$ws.Range("J2").Value = 1900
$ws.Range("K2").Value = 945
$ws.Range("L2").Value = 1900
$ws.Range("M2").Value = -832
$ws.Range("N2").Value = -2126
$ws.Range("H45").Value = 1531.8889
$ws.Range("I45").Value = 1766.6666
$ws.Range("J45").Value = 1484.9333
$ws.Range("K45").Value = 1766.6666
$ws.Range("L45").Value = 1484.9333
$ws.Range("M45").Value = -1389.6666
$ws.Range("N45").Value = -2238.9333
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H102").Value = 1448.8334
$ws.Range("I102").Value = 1398.7273
$ws.Range("K102").Value = 1398.7273
$ws.Range("M102").Value = 223.2727
$ws.Range("H110").Value = 1090.3334
$ws.Range("I110").Value = 1179.4
$ws.Range("J110").Value = 979
$ws.Range("K110").Value = 1179.4
$ws.Range("L110").Value = 979
$ws.Range("M110").Value = 865.5999999999999
$ws.Range("N110").Value = -5069
$ws.Range("H116").Value = 1250.6
$ws.Range("I116").Value = 945
$ws.Range("J116").Value = 1900
$ws.Range("K116").Value = 945
$ws.Range("L116").Value = 1900
$ws.Range("M116").Value = 1349
$ws.Range("N116").Value = -6488
$ws.Range("H132").Value = 2465.6428
$ws.Range("I132").Value = 2491.2112
$ws.Range("J132").Value = 2326
$ws.Range("K132").Value = 7473.633600000001
$ws.Range("L132").Value = 6978
$ws.Range("M132").Value = -4943.633600000001
$ws.Range("N132").Value = -12038

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1250.6
$ws.Range("I3").Value = 945
$ws.Range("J3").Value = 1900
$ws.Range("K3").Value = 945
$ws.Range("L3").Value = 1900
$ws.Range("M3").Value = -831
$ws.Range("N3").Value = -2128
$ws.Range("H105").Value = 2742.261
$ws.Range("I105").Value = 2389.0715
$ws.Range("J105").Value = 3291.6667
$ws.Range("K105").Value = 2389.0715
$ws.Range("L105").Value = 3291.6667
$ws.Range("M105").Value = -642.0715
$ws.Range("N105").Value = -6785.6667
$ws.Range("H107").Value = 748.43475
$ws.Range("I107").Value = 637.9286
$ws.Range("J107").Value = 920.3333
$ws.Range("K107").Value = 637.9286
$ws.Range("L107").Value = 920.3333
$ws.Range("M107").Value = 1282.0714
$ws.Range("N107").Value = -4760.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 915.5333000000001
$ws.Range("I16").Value = 827.75
$ws.Range("J16").Value = 1266.6666
$ws.Range("K16").Value = 827.75
$ws.Range("L16").Value = 1266.6666
$ws.Range("M16").Value = -540.75
$ws.Range("N16").Value = -1840.6666
$ws.Range("H22").Value = 465.47058
$ws.Range("I22").Value = 428.36365
$ws.Range("J22").Value = 533.5
$ws.Range("K22").Value = 428.36365
$ws.Range("L22").Value = 533.5
$ws.Range("M22").Value = -78.36365000000001
$ws.Range("N22").Value = -1233.5
$ws.Range("H107").Value = 642.6222
$ws.Range("I107").Value = 687.5833
$ws.Range("J107").Value = 591.2381
$ws.Range("K107").Value = 687.5833
$ws.Range("L107").Value = 591.2381
$ws.Range("M107").Value = 1232.4167
$ws.Range("N107").Value = -4431.2381
$ws.Range("H113").Value = 915.5333000000001
$ws.Range("I113").Value = 827.75
$ws.Range("J113").Value = 1266.6666
$ws.Range("K113").Value = 827.75
$ws.Range("L113").Value = 1266.6666
$ws.Range("M113").Value = 1342.25
$ws.Range("N113").Value = -5606.6666
$ws.Range("H134").Value = 1094.826
$ws.Range("I134").Value = 720.0526
$ws.Range("K134").Value = 2160.1578
$ws.Range("M134").Value = 374.8422

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 401.9524
$ws.Range("I2").Value = 442.1579
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 2652.9474
$ws.Range("L2").Value = 120
$ws.Range("M2").Value = -2539.9474
$ws.Range("N2").Value = -346
$ws.Range("H42").Value = 2799.5
$ws.Range("I42").Value = 2799
$ws.Range("J42").Value = 2800
$ws.Range("K42").Value = 8397
$ws.Range("L42").Value = 8400
$ws.Range("M42").Value = -7863
$ws.Range("N42").Value = -9468
$ws.Range("H114").Value = 5968.1904
$ws.Range("I114").Value = 11200.1
$ws.Range("J114").Value = 1211.909
$ws.Range("K114").Value = 33600.3
$ws.Range("L114").Value = 3635.727
$ws.Range("M114").Value = -30346.3
$ws.Range("N114").Value = -10143.727
$ws.Range("H134").Value = 1803.0769
$ws.Range("I134").Value = 1493.3334
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 4480.0002
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = 589.9997999999996
$ws.Range("N134").Value = -17640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 975.7368
$ws.Range("I97").Value = 861.5833
$ws.Range("K97").Value = 861.5833
$ws.Range("M97").Value = -365.5833
$ws.Range("H107").Value = 306.15
$ws.Range("J107").Value = 278.45456
$ws.Range("L107").Value = 278.45456
$ws.Range("N107").Value = -4118.45456
$ws.Range("H113").Value = 31250748
$ws.Range("I113").Value = 62500496
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 62500496
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = -62498326
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 41668716
$ws.Range("I7").Value = 1596.6666
$ws.Range("J7").Value = 83335830
$ws.Range("K7").Value = 1596.6666
$ws.Range("L7").Value = 83335830
$ws.Range("M7").Value = -1484.6666
$ws.Range("N7").Value = -83336054
$ws.Range("H61").Value = 10417812
$ws.Range("I61").Value = 1065.68
$ws.Range("J61").Value = 47620476
$ws.Range("K61").Value = 1065.68
$ws.Range("L61").Value = 47620476
$ws.Range("M61").Value = -863.6800000000001
$ws.Range("N61").Value = -47620880
$ws.Range("H113").Value = 10417812
$ws.Range("I113").Value = 1065.68
$ws.Range("J113").Value = 47620476
$ws.Range("K113").Value = 1065.68
$ws.Range("L113").Value = 47620476
$ws.Range("M113").Value = 1104.32
$ws.Range("N113").Value = -47624816
$ws.Range("H122").Value = 2177.0688
$ws.Range("I122").Value = 1863.5714
$ws.Range("K122").Value = 5590.7142
$ws.Range("M122").Value = -3140.7142
$ws.Range("H126").Value = 41668716
$ws.Range("I126").Value = 1596.6666
$ws.Range("J126").Value = 83335830
$ws.Range("K126").Value = 4789.9998
$ws.Range("L126").Value = 250007490
$ws.Range("M126").Value = -2319.9998
$ws.Range("N126").Value = -250012430

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 437.57895
$ws.Range("I107").Value = 507.36365
$ws.Range("J107").Value = 341.625
$ws.Range("K107").Value = 1522.09095
$ws.Range("L107").Value = 1024.875
$ws.Range("M107").Value = 397.90905
$ws.Range("N107").Value = -4864.875
$ws.Range("H113").Value = 563.5714
$ws.Range("I113").Value = 667.25
$ws.Range("J113").Value = 425.33334
$ws.Range("K113").Value = 2001.75
$ws.Range("L113").Value = 1276.00002
$ws.Range("M113").Value = 168.25
$ws.Range("N113").Value = -5616.000019999999
